$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is a new "classification" snapshot row, identical in shape to the
# existing rows 2/3 but with an updated timestamp and market status (and
# matching "未开盘" pre-market values). Copy row 3 down to row 4 first so the
# text-typed columns (A-F) keep their original text cell type/format, then
# overwrite the two cells that actually differ.
$ws.Range("A3:U3").Copy()
$ws.Range("A4").PasteSpecial()

$ws.Range("A4").Value = "Tue Oct 31 00:50:03 2023"
$ws.Range("D4").Value = "未开盘"
